# Rewrite the single " m: (...).fromTextileString(...) " field (fldChar
# begin / instrText* / fldChar end) that lives alone in its own paragraph
# into plain literal text runs "{ m: (...) }" (curly braces replacing the
# field delimiters), while keeping the _GoBack bookmark in the same
# relative position between the runs.

$d = $word.ActiveDocument

# Locate the field (there is exactly one in this document) and the
# paragraph that solely contains it.
$fld = $d.Fields.Item(1)
$fldCodeStart = $fld.Code.Start

$pos = -1
foreach ($p in $d.Paragraphs) {
    if (($fldCodeStart -ge $p.Range.Start) -and ($fldCodeStart -lt $p.Range.End)) {
        $pos = $p.Range.Start
    }
}

# Remove the whole field (fldChar begin, every instrText run, the
# bookmark and fldChar end) -- we rebuild the paragraph from scratch.
$fld.Delete()

# The text content of each resulting <w:r><w:t>...</w:t></w:r> run, in
# order.
$runTexts = @(
    '{',
    'm',
    ':',
    ' (',
    'self.na',
    'me',
    ' + ''\n\n!../images/logo_M2Doc.png',
    '!''',
    ').from',
    'Textile',
    'String(',
    '''https://www.m2doc.org/tests/''',
    ')',
    '}'
)

# Index (0-based, into $runTexts) after which the _GoBack bookmark must
# be re-inserted -- i.e. right after the ".png" run and before "!'".
$bookmarkAfterIndex = 6

# Pass 1: type every run out, separated by a uniquely-named temporary
# bookmark so each InsertAfter starts a brand new run instead of being
# merged into the previous one (adjacent same-format runs get coalesced
# otherwise). Marking the _GoBack split point with its own temp name too
# -- adding a bookmark exactly at the trailing edge of what has been
# typed so far (nothing after it yet) can make it snap to the wrong
# paragraph, so the permanent bookmark is only added in pass 2, once
# real content exists on both sides of it.
$tempBookmarks = New-Object System.Collections.ArrayList
$goBackPos = -1

for ($i = 0; $i -lt $runTexts.Length; $i++) {
    $ins = $d.Range($pos, $pos)
    $ins.InsertAfter($runTexts[$i])
    $pos = $pos + $runTexts[$i].Length

    if ($i -lt $runTexts.Length - 1) {
        $bmName = "TempRunSplit$i"
        $bmRange = $d.Range($pos, $pos)
        $d.Bookmarks.Add($bmName, $bmRange) | Out-Null
        [void]$tempBookmarks.Add($bmName)

        if ($i -eq $bookmarkAfterIndex) {
            $goBackPos = $pos
        }
    }
}

# Pass 2: drop every temporary split marker, then place the real
# _GoBack bookmark at the remembered position (now safely mid-paragraph).
foreach ($bmName in $tempBookmarks) {
    $d.Bookmarks.Item($bmName).Delete()
}

$goBackRange = $d.Range($goBackPos, $goBackPos)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null
